$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build list of target cells that need text updates (Price and Volume columns)
# Force text number format first so Excel does not auto-convert numeric-looking
# strings (e.g. "1.002", "0.9250", "0.000007972") into actual numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.701.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.920.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.26"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4938"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2981"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06782"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.905.97"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.22"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07352"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.174"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.96"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6738"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.678.66"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007972"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.62"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.142.97"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.337"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +10.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "200.69"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.321"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.666"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.72"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +7.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.94"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.968"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.53%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.379"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09182"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.067"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05302"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7445"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.121"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.733"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01842"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9250"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.92"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +29.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4466"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.971"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.02"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1390"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.650"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.95"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.90%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05879"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.58%  "
